$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.148.93"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").Value = "2.840.05"
$ws.Range("E3").Value = "  +1.87%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "362.04"
$ws.Range("E5").Value = "  +6.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.02"
$ws.Range("E6").Value = "  -2.65%  "

$ws.Range("E7").Value = "  +3.91%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.604"
$ws.Range("E9").Value = "  +4.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.52"
$ws.Range("E10").Value = "  -1.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0863"
$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("E12").Value = "  +1.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.03"
$ws.Range("E13").Value = "  -0.21%  "

$ws.Range("E14").Value = "  +2.58%  "

$ws.Range("D15").Value = "3.283.87"
$ws.Range("E15").Value = "  +1.66%  "

$ws.Range("D16").Value = "2.830.87"
$ws.Range("E16").Value = "  +0.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.910"
$ws.Range("E17").Value = "  +3.03%  "

$ws.Range("D18").Value = "52.062.27"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("E19").Value = "  +8.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.16"
$ws.Range("E20").Value = "  -1.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.56"
$ws.Range("E21").Value = "  +2.26%  "

$ws.Range("D22").Value = "0.0₃0998"
$ws.Range("E22").Value = "  +1.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.35"
$ws.Range("E23").Value = "  +0.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.67"
$ws.Range("E24").Value = "  -3.48%  "

$ws.Range("E25").Value = "  +3.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.09"
$ws.Range("E26").Value = "  +1.29%  "

$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("E28").Value = "  +1.87%  "

$ws.Range("E29").Value = "  +1.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0489"
$ws.Range("E30").Value = "  +30.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "53.95"
$ws.Range("E31").Value = "  +7.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.140"
$ws.Range("E32").Value = "  -1.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.30"
$ws.Range("E33").Value = "  +1.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.88"
$ws.Range("E34").Value = "  +2.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.49"
$ws.Range("E35").Value = "  +10.94%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0845"
$ws.Range("E36").Value = "  +2.63%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.27"
$ws.Range("E38").Value = "  +0.90%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.07"
$ws.Range("E39").Value = "  -2.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.45"
$ws.Range("E40").Value = "  -2.30%  "

$ws.Range("E41").Value = "  +1.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.66"
$ws.Range("E42").Value = "  +2.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.22"
$ws.Range("E43").Value = "  +1.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.53"
$ws.Range("E44").Value = "  -7.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.27"
$ws.Range("E45").Value = "  -3.40%  "

$ws.Range("E46").Value = "  +3.11%  "

$ws.Range("D47").Value = "2.113.75"
$ws.Range("E47").Value = "  +1.17%  "

$ws.Range("E48").Value = "  +1.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.997"
$ws.Range("E49").Value = "  +11.74%  "

$ws.Range("E50").Value = "  +5.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.25"
$ws.Range("E51").Value = "  +4.22%  "
